# Add "Extension Script file" column (with the Ansible WinRM configuration
# script link) to the Azure app invoice worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at C, shifting the existing data (OS name, images,
# sizes, regions, ...) one column to the right.
$ws.Columns("C:C").Insert()

$url = "https://github.com/ansible/ansible/blob/devel/examples/scripts/ConfigureRemotingForAnsible.ps1"

$ws.Range("C1").Value = "Extension Script file"
$ws.Range("C2").Value = $url
$ws.Range("C3").Value = $url

# Workbook-level bookkeeping that Excel touches when the file is re-saved
# from its new location.
$wb.Sheets.Item(1).Select()
$ws.Range("C22").Select()
